$d = $word.ActiveDocument

# --- document.xml: split the single "{code}" paragraph into two ---
# 1. Rewrite the leading run's text in place (keeps the bookmark + the
#    trailing "}" run untouched) so paragraph 1 becomes "{#list}{content}".
$null = $d.Content.Find.Execute("{code", $true, $false, $false, $false, $false, $true, 1, $false, "{#list}{content", 2)

# 2. Split right after that paragraph, producing a brand-new second
#    paragraph, then give it the "{/list}" text.
$p1 = $d.Paragraphs(1)
$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs(2)
$p2.Range.Text = "{/list}"

# --- styles.xml: mark "HTML Preformatted" (paragraph + linked character
#     style) as Quick Styles (adds <w:qFormat/>) ---
$d.Styles("HTML Preformatted").QuickStyle = $true
$d.Styles.Item(10).QuickStyle = $true
